# Actualización automática hashcode jue ene 17 01:35:07 CET 2019
# Updates the hashcode values (column B) for specific rows identified by
# their previous hash value, replacing them with newly computed hashes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = "B44";  New = "801441973795329981791b8136476d51" },
    @{ Cell = "B51";  New = "06f8099ed22fd4689878c01d8241bac2" },
    @{ Cell = "B74";  New = "81b8198663d8342ceb3b8c0f92fab114" },
    @{ Cell = "B80";  New = "7cd71806c3817a2788b411cc5dc0d07f" },
    @{ Cell = "B89";  New = "1616b5e7f8bed5b4d7aed86321c8e87e" },
    @{ Cell = "B99";  New = "934acdaaaa0b3be31f1a4c83585356c0" },
    @{ Cell = "B108"; New = "c837468acc659d7ed0d988fd25708386" },
    @{ Cell = "B110"; New = "391e31b1a8bd2400f63b4fbdf2ed30bd" },
    @{ Cell = "B121"; New = "27ce3918723a74c22be7d3b4776af7d0" },
    @{ Cell = "B161"; New = "43b27c02768b9c7c3fa9e56208ca190b" },
    @{ Cell = "B168"; New = "a1b0e2550e24d1d6623b2a13cb8c46cb" },
    @{ Cell = "B278"; New = "c471259a9ae3506bba77c0b291834b56" },
    @{ Cell = "B345"; New = "d1f32890b74c9e8aba42588b693f86cc" },
    @{ Cell = "B540"; New = "99f4011882d24541f61623d2c1b5defc" },
    @{ Cell = "B543"; New = "ef5f9019c2a4a7b02d1df030ca1ce0aa" },
    @{ Cell = "B574"; New = "a129a870088d76f781fe1f5950d3a8ba" },
    @{ Cell = "B575"; New = "03f38022c575245c28fc04992de3c384" },
    @{ Cell = "B616"; New = "cf51451dd6f5b3073cd680b0a9c8f098" },
    @{ Cell = "B715"; New = "d6ec5b2a28c05cafb949242c8f5515d0" },
    @{ Cell = "B768"; New = "b45c8bde2cac9396d620eb045d985164" },
    @{ Cell = "B816"; New = "dc3ff660a48a009b2c263afaeeb131db" },
    @{ Cell = "B825"; New = "ee144aaf330dcd969107a5068c1f5d28" },
    @{ Cell = "B827"; New = "b12f29376da282e56a56ae942e4a5f02" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.New
}
